$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H31").Value = 634
$ws.Range("I31").Value = 634
$ws.Range("K31").Value = 1902
$ws.Range("M31").Value = -1672

$ws.Range("H32").Value = 538.5
$ws.Range("I32").Value = 800
$ws.Range("J32").Value = 419.63635
$ws.Range("K32").Value = 800
$ws.Range("L32").Value = 419.63635
$ws.Range("M32").Value = -474
$ws.Range("N32").Value = -1071.63635

$ws.Range("H137").Value = 11364845
$ws.Range("I137").Value = 15152274
$ws.Range("J137").Value = 2558.3635
$ws.Range("K137").Value = 45456822
$ws.Range("L137").Value = 7675.0905
$ws.Range("M137").Value = -45454272
$ws.Range("N137").Value = -12775.0905

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 104.625
$ws.Range("I4").Value = 49.4
$ws.Range("J4").Value = 196.66667
$ws.Range("K4").Value = 49.4
$ws.Range("L4").Value = 196.66667
$ws.Range("M4").Value = 66.59999999999999
$ws.Range("N4").Value = -428.66667

$ws.Range("H5").Value = 9523866
$ws.Range("I5").Value = 12987059
$ws.Range("J5").Value = 84.5
$ws.Range("K5").Value = 12987059
$ws.Range("L5").Value = 84.5
$ws.Range("M5").Value = -12986947
$ws.Range("N5").Value = -308.5

$ws.Range("H8").Value = 10000
$ws.Range("I8").Value = 0
$ws.Range("J8").Value = 10000
$ws.Range("K8").Value = 0
$ws.Range("L8").Value = 10000
$ws.Range("M8").Value = ""
$ws.Range("N8").Value = -10288

$ws.Range("H31").Value = 10389.333
$ws.Range("I31").Value = 1985.7778
$ws.Range("J31").Value = 35600
$ws.Range("K31").Value = 1985.7778
$ws.Range("L31").Value = 35600
$ws.Range("M31").Value = -1691.7778
$ws.Range("N31").Value = -36188

$ws.Range("H39").Value = 7750
$ws.Range("I39").Value = 8000
$ws.Range("J39").Value = 7500
$ws.Range("K39").Value = 8000
$ws.Range("L39").Value = 7500
$ws.Range("M39").Value = -7480
$ws.Range("N39").Value = -8540

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 9523866
$ws.Range("I4").Value = 12987059
$ws.Range("J4").Value = 84.5
$ws.Range("K4").Value = 12987059
$ws.Range("L4").Value = 84.5
$ws.Range("M4").Value = -12986944
$ws.Range("N4").Value = -314.5

$ws.Range("H21").Value = 21833.334
$ws.Range("J21").Value = 21833.334
$ws.Range("L21").Value = 21833.334
$ws.Range("N21").Value = -22305.334

$ws.Range("H56").Value = 46833.332
$ws.Range("J56").Value = 46833.332
$ws.Range("L56").Value = 46833.332
$ws.Range("N56").Value = -48311.332

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 64.25
$ws.Range("I7").Value = 87.55556
$ws.Range("J7").Value = 34.285713
$ws.Range("K7").Value = 87.55556
$ws.Range("L7").Value = 34.285713
$ws.Range("M7").Value = 25.44444
$ws.Range("N7").Value = -260.285713

$ws.Range("H22").Value = 1000.125
$ws.Range("I22").Value = 600.5
$ws.Range("J22").Value = 1133.3334
$ws.Range("K22").Value = 600.5
$ws.Range("L22").Value = 1133.3334
$ws.Range("M22").Value = -250.5
$ws.Range("N22").Value = -1833.3334

$ws.Range("H36").Value = 37426.5
$ws.Range("I36").Value = 0
$ws.Range("J36").Value = 37426.5
$ws.Range("K36").Value = 0
$ws.Range("L36").Value = 37426.5
$ws.Range("M36").Value = ""
$ws.Range("N36").Value = -38202.5

$ws.Range("H40").Value = 37426.5
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 37426.5
$ws.Range("K40").Value = 0
$ws.Range("L40").Value = 37426.5
$ws.Range("M40").Value = ""
$ws.Range("N40").Value = -37746.5

$ws.Range("H100").Value = 38943.6
$ws.Range("J100").Value = 38943.6
$ws.Range("L100").Value = 38943.6
$ws.Range("N100").Value = -41107.6

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H42").Value = 0
$ws.Range("I42").Value = 0
$ws.Range("J42").Value = 0
$ws.Range("K42").Value = 0
$ws.Range("L42").Value = 0
$ws.Range("M42").Value = ""
$ws.Range("N42").Value = ""

$ws.Range("H48").Value = 0
$ws.Range("J48").Value = 0
$ws.Range("L48").Value = 0
$ws.Range("N48").Value = ""

$ws.Range("H114").Value = 1722.6154
$ws.Range("I114").Value = 1471
$ws.Range("J114").Value = 1834.4445
$ws.Range("K114").Value = 4413
$ws.Range("L114").Value = 5503.333500000001
$ws.Range("M114").Value = -1159
$ws.Range("N114").Value = -12011.3335

$ws.Range("H117").Value = 2997
$ws.Range("J117").Value = 3796
$ws.Range("L117").Value = 11388
$ws.Range("N117").Value = -18272

$ws.Range("H121").Value = 52637576
$ws.Range("I121").Value = 610
$ws.Range("J121").Value = 71436500
$ws.Range("K121").Value = 1830
$ws.Range("L121").Value = 214309500
$ws.Range("M121").Value = -520
$ws.Range("N121").Value = -214312120

$ws.Range("H136").Value = 779.3
$ws.Range("I136").Value = 570
$ws.Range("K136").Value = 1710
$ws.Range("M136").Value = 3390

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 785.3611
$ws.Range("I107").Value = 717.08
$ws.Range("K107").Value = 717.08
$ws.Range("M107").Value = 1202.92

$ws.Range("H126").Value = 1964.7273
$ws.Range("I126").Value = 1608.2667
$ws.Range("J126").Value = 2728.5715
$ws.Range("K126").Value = 4824.800099999999
$ws.Range("L126").Value = 8185.7145
$ws.Range("M126").Value = -2354.800099999999
$ws.Range("N126").Value = -13125.7145

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 942
$ws.Range("I22").Value = 973.3333
$ws.Range("J22").Value = 928.5714
$ws.Range("K22").Value = 973.3333
$ws.Range("L22").Value = 928.5714
$ws.Range("M22").Value = -678.3333
$ws.Range("N22").Value = -1518.5714

$ws.Range("H27").Value = 942
$ws.Range("I27").Value = 973.3333
$ws.Range("J27").Value = 928.5714
$ws.Range("K27").Value = 973.3333
$ws.Range("L27").Value = 928.5714
$ws.Range("M27").Value = -866.3333
$ws.Range("N27").Value = -1142.5714

$ws.Range("H33").Value = 5708.5
$ws.Range("I33").Value = 1400
$ws.Range("J33").Value = 10017
$ws.Range("K33").Value = 1400
$ws.Range("L33").Value = 10017
$ws.Range("M33").Value = -1110
$ws.Range("N33").Value = -10597

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H51").Value = 10000
$ws.Range("I51").Value = 10000
$ws.Range("K51").Value = 10000
$ws.Range("M51").Value = -9490

$ws.Range("H62").Value = 2579.3125
$ws.Range("I62").Value = 1915
$ws.Range("J62").Value = 2977.9
$ws.Range("K62").Value = 1915
$ws.Range("L62").Value = 2977.9
$ws.Range("M62").Value = -1291
$ws.Range("N62").Value = -4225.9

$ws.Range("H65").Value = 2579.3125
$ws.Range("I65").Value = 1915
$ws.Range("J65").Value = 2977.9
$ws.Range("K65").Value = 9575
$ws.Range("L65").Value = 14889.5
$ws.Range("M65").Value = -6455
$ws.Range("N65").Value = -21129.5
